$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "datos actualizados" timestamp title (A1)
$ws.Range("A1").Value = "Datos actualizados a 14 de Agosto de 2020 a las 01:16"

# Country name swaps caused by re-sorting adjacent countries with near-identical totals
$ws.Range("A106").Value = "Zimbabue"
$ws.Range("A107").Value = "Malaui"
$ws.Range("A115").Value = "Suazilandia"
$ws.Range("A116").Value = "Namibia"
$ws.Range("A123").Value = "Surinam"
$ws.Range("A124").Value = "Eslovaquia"
$ws.Range("A169").Value = "Trinidad yTobago"
$ws.Range("A170").Value = "Comoras"
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("A214").Value = "Montserrat"

# Update daily statistics (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
$ws.Range("B4").Value = 5410938
$ws.Range("C4").Value = 50636
$ws.Range("D4").Value = 2832524
$ws.Range("E4").Value = 2408155
$ws.Range("G4").Value = 1128
$ws.Range("H4").Value = 170259
$ws.Range("B11").Value = 433805
$ws.Range("C11").Value = 11286
$ws.Range("D11").Value = 250494
$ws.Range("E11").Value = 169166
$ws.Range("G11").Value = 308
$ws.Range("H11").Value = 14145
$ws.Range("B31").Value = 96108
$ws.Range("C31").Value = 145
$ws.Range("D31").Value = 56890
$ws.Range("E31").Value = 34111
$ws.Range("G31").Value = 22
$ws.Range("H31").Value = 5107
$ws.Range("B46").Value = 60284
$ws.Range("C46").Value = 1195
$ws.Range("D46").Value = 48305
$ws.Range("E46").Value = 9683
$ws.Range("G46").Value = 29
$ws.Range("H46").Value = 2296
$ws.Range("B50").Value = 51147
$ws.Range("C50").Value = 937
$ws.Range("D50").Value = 36134
$ws.Range("E50").Value = 13950
$ws.Range("G50").Value = 4
$ws.Range("H50").Value = 1063
$ws.Range("B52").Value = 48116
$ws.Range("C52").Value = 373
$ws.Range("D52").Value = 34309
$ws.Range("E52").Value = 12841
$ws.Range("G52").Value = 10
$ws.Range("H52").Value = 966
$ws.Range("B74").Value = 19401
$ws.Range("C74").Value = 326
$ws.Range("D74").Value = 13574
$ws.Range("E74").Value = 5436
$ws.Range("B86").Value = 9851
$ws.Range("C86").Value = 68
$ws.Range("E86").Value = 737
$ws.Range("B100").Value = 6653
$ws.Range("C100").Value = 31
$ws.Range("D100").Value = 5843
$ws.Range("E100").Value = 653
$ws.Range("B106").Value = 4990
$ws.Range("C106").Value = 97
$ws.Range("D106").Value = 1927
$ws.Range("E106").Value = 2935
$ws.Range("G106").Value = 6
$ws.Range("H106").Value = 128
$ws.Range("B107").Value = 4912
$ws.Range("C107").Value = 160
$ws.Range("D107").Value = 2550
$ws.Range("E107").Value = 2209
$ws.Range("G107").Value = 1
$ws.Range("H107").Value = 153
$ws.Range("B113").Value = 3857
$ws.Range("C113").Value = 44
$ws.Range("D113").Value = 2680
$ws.Range("E113").Value = 1104
$ws.Range("B115").Value = 3599
$ws.Range("C115").Value = 74
$ws.Range("D115").Value = 1991
$ws.Range("E115").Value = 1543
$ws.Range("G115").Value = 2
$ws.Range("H115").Value = 65
$ws.Range("B116").Value = 3544
$ws.Range("C116").Value = 138
$ws.Range("D116").Value = 848
$ws.Range("E116").Value = 2669
$ws.Range("G116").Value = 5
$ws.Range("H116").Value = 27
$ws.Range("B123").Value = 2761
$ws.Range("C123").Value = 108
$ws.Range("D123").Value = 1830
$ws.Range("E123").Value = 891
$ws.Range("G123").Value = 1
$ws.Range("H123").Value = 40
$ws.Range("B124").Value = 2739
$ws.Range("C124").Value = 49
$ws.Range("D124").Value = 1939
$ws.Range("E124").Value = 769
$ws.Range("H124").Value = 31
$ws.Range("B127").Value = 2478
$ws.Range("C127").Value = 1
$ws.Range("E127").Value = 1256
$ws.Range("B158").Value = 911
$ws.Range("C158").Value = 28
$ws.Range("D158").Value = 425
$ws.Range("E158").Value = 465
$ws.Range("G158").Value = 4
$ws.Range("H158").Value = 21
$ws.Range("B165").Value = 631
$ws.Range("C165").Value = 8
$ws.Range("D165").Value = 202
$ws.Range("E165").Value = 407
$ws.Range("B169").Value = 404
$ws.Range("C169").Value = 78
$ws.Range("D169").Value = 139
$ws.Range("E169").Value = 257
$ws.Range("H169").Value = 8
$ws.Range("B170").Value = 399
$ws.Range("C170").Value = 0
$ws.Range("D170").Value = 379
$ws.Range("E170").Value = 13
$ws.Range("H170").Value = 7
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1
